$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-05 Tuesday", "2024-11-06 Wednesday"),
    @("414÷8=", "204÷2="),
    @("257÷6=", "485÷4="),
    @("434÷3=", "572÷8="),
    @("701÷6=", "494÷7="),
    @("628÷4=", "650÷2="),
    @("799÷5=", "820÷4="),
    @("537÷3=", "800÷9="),
    @("105÷9=", "441÷7="),
    @("476÷4=", "308÷4="),
    @("879÷6=", "854÷2="),
    @("960÷2=", "221÷9="),
    @("502÷7=", "407÷7="),
    @("811÷9=", "886÷9="),
    @("770÷8=", "641÷4="),
    @("733÷3=", "984÷5="),
    @("835÷6=", "360÷2="),
    @("346÷9=", "230÷5="),
    @("911÷7=", "650÷7="),
    @("726÷5=", "900÷3="),
    @("238÷4=", "921÷4="),
    @("365÷2=", "276÷3="),
    @("600÷8=", "818÷8="),
    @("647÷8=", "942÷7="),
    @("789÷9=", "742÷8="),
    @("137÷8=", "913÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
